$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 672.5
$ws.Range("I11").Value = 672.5
$ws.Range("K11").Value = 672.5
$ws.Range("M11").Value = -532.5

$ws.Range("H21").Value = 10019
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H23").Value = 10019
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H33").Value = 287.3158
$ws.Range("I33").Value = 247.72223
$ws.Range("K33").Value = 247.72223
$ws.Range("M33").Value = -18.72223

$ws.Range("H92").Value = 7617
$ws.Range("I92").Value = 5999
$ws.Range("J92").Value = 8587.799999999999
$ws.Range("K92").Value = 5999
$ws.Range("L92").Value = 8587.799999999999
$ws.Range("M92").Value = -4751
$ws.Range("N92").Value = -11083.8

$ws.Range("H137").Value = 3675.4443
$ws.Range("I137").Value = 2020.2
$ws.Range("J137").Value = 5744.5
$ws.Range("K137").Value = 6060.6
$ws.Range("L137").Value = 17233.5
$ws.Range("M137").Value = -3510.6
$ws.Range("N137").Value = -22333.5

$ws.Range("H141").Value = 6204
$ws.Range("I141").Value = 2010
$ws.Range("K141").Value = 6030
$ws.Range("M141").Value = -850

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1293.5103
$ws.Range("I2").Value = 1257.9791
$ws.Range("J2").Value = 2999
$ws.Range("K2").Value = 1257.9791
$ws.Range("L2").Value = 2999
$ws.Range("M2").Value = -1144.9791
$ws.Range("N2").Value = -3225

$ws.Range("H45").Value = 3978
$ws.Range("I45").Value = 3544.2856
$ws.Range("J45").Value = 7014
$ws.Range("K45").Value = 3544.2856
$ws.Range("L45").Value = 7014
$ws.Range("M45").Value = -3167.2856
$ws.Range("N45").Value = -7768

$ws.Range("H74").Value = 1680.25
$ws.Range("I74").Value = 1606
$ws.Range("J74").Value = 2200
$ws.Range("K74").Value = 1606
$ws.Range("L74").Value = 2200
$ws.Range("M74").Value = -732
$ws.Range("N74").Value = -3948

$ws.Range("H77").Value = 1680.25
$ws.Range("I77").Value = 1606
$ws.Range("J77").Value = 2200
$ws.Range("K77").Value = 8030
$ws.Range("L77").Value = 11000
$ws.Range("M77").Value = -3662
$ws.Range("N77").Value = -19736

$ws.Range("H116").Value = 1293.5103
$ws.Range("I116").Value = 1257.9791
$ws.Range("J116").Value = 2999
$ws.Range("K116").Value = 1257.9791
$ws.Range("L116").Value = 2999
$ws.Range("M116").Value = 1036.0209
$ws.Range("N116").Value = -7587

$ws.Range("H122").Value = 2569.2354
$ws.Range("J122").Value = 3100
$ws.Range("L122").Value = 9300
$ws.Range("N122").Value = -14200

$ws.Range("H132").Value = 2878.6052
$ws.Range("I132").Value = 2867.3635
$ws.Range("K132").Value = 8602.0905
$ws.Range("M132").Value = -6072.0905

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1293.5103
$ws.Range("I3").Value = 1257.9791
$ws.Range("J3").Value = 2999
$ws.Range("K3").Value = 1257.9791
$ws.Range("L3").Value = 2999
$ws.Range("M3").Value = -1143.9791
$ws.Range("N3").Value = -3227

$ws.Range("H105").Value = 3960.5652
$ws.Range("I105").Value = 3754.25
$ws.Range("J105").Value = 4004
$ws.Range("K105").Value = 3754.25
$ws.Range("L105").Value = 4004
$ws.Range("M105").Value = -2007.25
$ws.Range("N105").Value = -7498

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 6140.5713
$ws.Range("I2").Value = 5500
$ws.Range("K2").Value = 5500
$ws.Range("M2").Value = -5387

$ws.Range("H16").Value = 782.2857
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws.Range("H32").Value = 10010
$ws.Range("I32").Value = 10010
$ws.Range("K32").Value = 10010
$ws.Range("M32").Value = -9694

$ws.Range("H113").Value = 782.2857
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws.Range("H134").Value = 2539.4707
$ws.Range("I134").Value = 1989.5333
$ws.Range("J134").Value = 6664
$ws.Range("K134").Value = 5968.5999
$ws.Range("L134").Value = 19992
$ws.Range("M134").Value = -3433.5999
$ws.Range("N134").Value = -25062

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 914.35
$ws.Range("I34").Value = 369.88235
$ws.Range("J34").Value = 3999.6667
$ws.Range("K34").Value = 1109.64705
$ws.Range("L34").Value = 11999.0001
$ws.Range("M34").Value = -1025.64705
$ws.Range("N34").Value = -12167.0001

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H82").Value = 8625
$ws.Range("J82").Value = 8625
$ws.Range("L82").Value = 25875
$ws.Range("N82").Value = -26687

$ws.Range("H85").Value = 8625
$ws.Range("J85").Value = 8625
$ws.Range("L85").Value = 25875
$ws.Range("N85").Value = -28683

$ws.Range("H113").Value = 850
$ws.Range("I113").Value = 850
$ws.Range("K113").Value = 2550
$ws.Range("M113").Value = -380

$ws.Range("H122").Value = 2300.875
$ws.Range("J122").Value = 2168.3333
$ws.Range("L122").Value = 19514.9997
$ws.Range("N122").Value = -24414.9997

$ws.Range("H131").Value = 25381.717
$ws.Range("J131").Value = 1790.88
$ws.Range("L131").Value = 5372.64
$ws.Range("N131").Value = -15452.64

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9160.639999999999
$ws.Range("I80").Value = 5415.385
$ws.Range("K80").Value = 5415.385
$ws.Range("M80").Value = -4417.385

$ws.Range("H83").Value = 9160.639999999999
$ws.Range("I83").Value = 5415.385
$ws.Range("K83").Value = 27076.925
$ws.Range("M83").Value = -22084.925

$ws.Range("H102").Value = 58463.188
$ws.Range("I102").Value = 66029.42999999999
$ws.Range("K102").Value = 66029.42999999999
$ws.Range("M102").Value = -64407.42999999999

$ws.Range("H132").Value = 5365.8887
$ws.Range("I132").Value = 5365.8887
$ws.Range("K132").Value = 16097.6661
$ws.Range("M132").Value = -13567.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4164.048
$ws.Range("I40").Value = 4085.9443
$ws.Range("K40").Value = 4085.9443
$ws.Range("M40").Value = -3949.9443

$ws.Range("H82").Value = 2007.2222
$ws.Range("I82").Value = 1852.8462
$ws.Range("J82").Value = 2408.6
$ws.Range("K82").Value = 1852.8462
$ws.Range("L82").Value = 2408.6
$ws.Range("M82").Value = -1491.8462
$ws.Range("N82").Value = -3130.6

$ws.Range("H85").Value = 2007.2222
$ws.Range("I85").Value = 1852.8462
$ws.Range("J85").Value = 2408.6
$ws.Range("K85").Value = 1852.8462
$ws.Range("L85").Value = 2408.6
$ws.Range("M85").Value = -604.8462
$ws.Range("N85").Value = -4904.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 29249.75
$ws.Range("J63").Value = 48499.5
$ws.Range("L63").Value = 48499.5
$ws.Range("N63").Value = -49747.5

$ws.Range("H66").Value = 29249.75
$ws.Range("J66").Value = 48499.5
$ws.Range("L66").Value = 145498.5
$ws.Range("N66").Value = -151738.5

$ws.Range("H107").Value = 1324.5834
$ws.Range("J107").Value = 3869
$ws.Range("L107").Value = 11607
$ws.Range("N107").Value = -15447

$ws.Range("H136").Value = 3114.5
$ws.Range("I136").Value = 2970.0715
$ws.Range("J136").Value = 3620
$ws.Range("K136").Value = 8910.2145
$ws.Range("L136").Value = 10860
$ws.Range("M136").Value = -6360.2145
$ws.Range("N136").Value = -15960
